$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad" / changed-date) for every existing data row (2..426)
#    moves from 45192 to 45202 (serial date 2023-09-23 -> 2023-10-03).
for ($r = 2; $r -le 426; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2) Row 426 picks up an explicit row height (ht="15" customHeight="1") in the
#    target file - touch RowHeight so it is written out explicitly.
$ws.Rows.Item(426).RowHeight = 15

# 3) A brand-new row 427 is appended with a new permit entry.
$ws.Cells.Item(427, 1).Value = "A 46247-2023"

$ws.Cells.Item(427, 2).Value = 45196
$ws.Cells.Item(427, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(427, 3).Value = 45202
$ws.Cells.Item(427, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(427, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(427, 5).Value = "OSBY"

$ws.Cells.Item(427, 7).Value = 3.4
$ws.Cells.Item(427, 8).Value = 0
$ws.Cells.Item(427, 9).Value = 0
$ws.Cells.Item(427, 10).Value = 0
$ws.Cells.Item(427, 11).Value = 0
$ws.Cells.Item(427, 12).Value = 0
$ws.Cells.Item(427, 13).Value = 0
$ws.Cells.Item(427, 14).Value = 0
$ws.Cells.Item(427, 15).Value = 0
$ws.Cells.Item(427, 16).Value = 0
$ws.Cells.Item(427, 17).Value = 0

# Column R keeps the wrap-text formatting used throughout the sheet even
# though this row has no species text.
$ws.Cells.Item(427, 18).WrapText = $true
